$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix/realign fixture rows whose Home/Away data got shuffled: swap the full
# row payload (columns B:AB) between each pair of rows, keeping column A (the
# running match index) fixed to its own row.

# --- swap row 91 <-> row 92 ---
$B91_old = $ws.Range("B91").Value2
$C91_old = $ws.Range("C91").Value2
$D91_old = $ws.Range("D91").Value2
$E91_old = $ws.Range("E91").Value2
$F91_old = $ws.Range("F91").Value2
$G91_old = $ws.Range("G91").Value2
$H91_old = $ws.Range("H91").Value2
$I91_old = $ws.Range("I91").Value2
$J91_old = $ws.Range("J91").Value2
$K91_old = $ws.Range("K91").Value2
$L91_old = $ws.Range("L91").Value2
$M91_old = $ws.Range("M91").Value2
$N91_old = $ws.Range("N91").Value2
$O91_old = $ws.Range("O91").Value2
$P91_old = $ws.Range("P91").Value2
$Q91_old = $ws.Range("Q91").Value2
$R91_old = $ws.Range("R91").Value2
$S91_old = $ws.Range("S91").Value2
$T91_old = $ws.Range("T91").Value2
$U91_old = $ws.Range("U91").Value2
$V91_old = $ws.Range("V91").Value2
$W91_old = $ws.Range("W91").Value2
$X91_old = $ws.Range("X91").Value2
$Y91_old = $ws.Range("Y91").Value2
$Z91_old = $ws.Range("Z91").Value2
$AA91_old = $ws.Range("AA91").Value2
$AB91_old = $ws.Range("AB91").Value2
$B92_old = $ws.Range("B92").Value2
$C92_old = $ws.Range("C92").Value2
$D92_old = $ws.Range("D92").Value2
$E92_old = $ws.Range("E92").Value2
$F92_old = $ws.Range("F92").Value2
$G92_old = $ws.Range("G92").Value2
$H92_old = $ws.Range("H92").Value2
$I92_old = $ws.Range("I92").Value2
$J92_old = $ws.Range("J92").Value2
$K92_old = $ws.Range("K92").Value2
$L92_old = $ws.Range("L92").Value2
$M92_old = $ws.Range("M92").Value2
$N92_old = $ws.Range("N92").Value2
$O92_old = $ws.Range("O92").Value2
$P92_old = $ws.Range("P92").Value2
$Q92_old = $ws.Range("Q92").Value2
$R92_old = $ws.Range("R92").Value2
$S92_old = $ws.Range("S92").Value2
$T92_old = $ws.Range("T92").Value2
$U92_old = $ws.Range("U92").Value2
$V92_old = $ws.Range("V92").Value2
$W92_old = $ws.Range("W92").Value2
$X92_old = $ws.Range("X92").Value2
$Y92_old = $ws.Range("Y92").Value2
$Z92_old = $ws.Range("Z92").Value2
$AA92_old = $ws.Range("AA92").Value2
$AB92_old = $ws.Range("AB92").Value2
$ws.Range("B91").Value2 = $B92_old
$ws.Range("C91").Value2 = $C92_old
$ws.Range("D91").Value2 = $D92_old
$ws.Range("E91").Value2 = $E92_old
$ws.Range("F91").Value2 = $F92_old
$ws.Range("G91").Value2 = $G92_old
$ws.Range("H91").Value2 = $H92_old
$ws.Range("I91").Value2 = $I92_old
$ws.Range("J91").Value2 = $J92_old
$ws.Range("K91").Value2 = $K92_old
$ws.Range("L91").Value2 = $L92_old
$ws.Range("M91").Value2 = $M92_old
$ws.Range("N91").Value2 = $N92_old
$ws.Range("O91").Value2 = $O92_old
$ws.Range("P91").Value2 = $P92_old
$ws.Range("Q91").Value2 = $Q92_old
$ws.Range("R91").Value2 = $R92_old
$ws.Range("S91").Value2 = $S92_old
$ws.Range("T91").Value2 = $T92_old
$ws.Range("U91").Value2 = $U92_old
$ws.Range("V91").Value2 = $V92_old
$ws.Range("W91").Value2 = $W92_old
$ws.Range("X91").Value2 = $X92_old
$ws.Range("Y91").Value2 = $Y92_old
$ws.Range("Z91").Value2 = $Z92_old
$ws.Range("AA91").Value2 = $AA92_old
$ws.Range("AB91").Value2 = $AB92_old
$ws.Range("B92").Value2 = $B91_old
$ws.Range("C92").Value2 = $C91_old
$ws.Range("D92").Value2 = $D91_old
$ws.Range("E92").Value2 = $E91_old
$ws.Range("F92").Value2 = $F91_old
$ws.Range("G92").Value2 = $G91_old
$ws.Range("H92").Value2 = $H91_old
$ws.Range("I92").Value2 = $I91_old
$ws.Range("J92").Value2 = $J91_old
$ws.Range("K92").Value2 = $K91_old
$ws.Range("L92").Value2 = $L91_old
$ws.Range("M92").Value2 = $M91_old
$ws.Range("N92").Value2 = $N91_old
$ws.Range("O92").Value2 = $O91_old
$ws.Range("P92").Value2 = $P91_old
$ws.Range("Q92").Value2 = $Q91_old
$ws.Range("R92").Value2 = $R91_old
$ws.Range("S92").Value2 = $S91_old
$ws.Range("T92").Value2 = $T91_old
$ws.Range("U92").Value2 = $U91_old
$ws.Range("V92").Value2 = $V91_old
$ws.Range("W92").Value2 = $W91_old
$ws.Range("X92").Value2 = $X91_old
$ws.Range("Y92").Value2 = $Y91_old
$ws.Range("Z92").Value2 = $Z91_old
$ws.Range("AA92").Value2 = $AA91_old
$ws.Range("AB92").Value2 = $AB91_old

# --- swap row 110 <-> row 111 ---
$B110_old = $ws.Range("B110").Value2
$C110_old = $ws.Range("C110").Value2
$D110_old = $ws.Range("D110").Value2
$E110_old = $ws.Range("E110").Value2
$F110_old = $ws.Range("F110").Value2
$G110_old = $ws.Range("G110").Value2
$H110_old = $ws.Range("H110").Value2
$I110_old = $ws.Range("I110").Value2
$J110_old = $ws.Range("J110").Value2
$K110_old = $ws.Range("K110").Value2
$L110_old = $ws.Range("L110").Value2
$M110_old = $ws.Range("M110").Value2
$N110_old = $ws.Range("N110").Value2
$O110_old = $ws.Range("O110").Value2
$P110_old = $ws.Range("P110").Value2
$Q110_old = $ws.Range("Q110").Value2
$R110_old = $ws.Range("R110").Value2
$S110_old = $ws.Range("S110").Value2
$T110_old = $ws.Range("T110").Value2
$U110_old = $ws.Range("U110").Value2
$V110_old = $ws.Range("V110").Value2
$W110_old = $ws.Range("W110").Value2
$X110_old = $ws.Range("X110").Value2
$Y110_old = $ws.Range("Y110").Value2
$Z110_old = $ws.Range("Z110").Value2
$AA110_old = $ws.Range("AA110").Value2
$AB110_old = $ws.Range("AB110").Value2
$B111_old = $ws.Range("B111").Value2
$C111_old = $ws.Range("C111").Value2
$D111_old = $ws.Range("D111").Value2
$E111_old = $ws.Range("E111").Value2
$F111_old = $ws.Range("F111").Value2
$G111_old = $ws.Range("G111").Value2
$H111_old = $ws.Range("H111").Value2
$I111_old = $ws.Range("I111").Value2
$J111_old = $ws.Range("J111").Value2
$K111_old = $ws.Range("K111").Value2
$L111_old = $ws.Range("L111").Value2
$M111_old = $ws.Range("M111").Value2
$N111_old = $ws.Range("N111").Value2
$O111_old = $ws.Range("O111").Value2
$P111_old = $ws.Range("P111").Value2
$Q111_old = $ws.Range("Q111").Value2
$R111_old = $ws.Range("R111").Value2
$S111_old = $ws.Range("S111").Value2
$T111_old = $ws.Range("T111").Value2
$U111_old = $ws.Range("U111").Value2
$V111_old = $ws.Range("V111").Value2
$W111_old = $ws.Range("W111").Value2
$X111_old = $ws.Range("X111").Value2
$Y111_old = $ws.Range("Y111").Value2
$Z111_old = $ws.Range("Z111").Value2
$AA111_old = $ws.Range("AA111").Value2
$AB111_old = $ws.Range("AB111").Value2
$ws.Range("B110").Value2 = $B111_old
$ws.Range("C110").Value2 = $C111_old
$ws.Range("D110").Value2 = $D111_old
$ws.Range("E110").Value2 = $E111_old
$ws.Range("F110").Value2 = $F111_old
$ws.Range("G110").Value2 = $G111_old
$ws.Range("H110").Value2 = $H111_old
$ws.Range("I110").Value2 = $I111_old
$ws.Range("J110").Value2 = $J111_old
$ws.Range("K110").Value2 = $K111_old
$ws.Range("L110").Value2 = $L111_old
$ws.Range("M110").Value2 = $M111_old
$ws.Range("N110").Value2 = $N111_old
$ws.Range("O110").Value2 = $O111_old
$ws.Range("P110").Value2 = $P111_old
$ws.Range("Q110").Value2 = $Q111_old
$ws.Range("R110").Value2 = $R111_old
$ws.Range("S110").Value2 = $S111_old
$ws.Range("T110").Value2 = $T111_old
$ws.Range("U110").Value2 = $U111_old
$ws.Range("V110").Value2 = $V111_old
$ws.Range("W110").Value2 = $W111_old
$ws.Range("X110").Value2 = $X111_old
$ws.Range("Y110").Value2 = $Y111_old
$ws.Range("Z110").Value2 = $Z111_old
$ws.Range("AA110").Value2 = $AA111_old
$ws.Range("AB110").Value2 = $AB111_old
$ws.Range("B111").Value2 = $B110_old
$ws.Range("C111").Value2 = $C110_old
$ws.Range("D111").Value2 = $D110_old
$ws.Range("E111").Value2 = $E110_old
$ws.Range("F111").Value2 = $F110_old
$ws.Range("G111").Value2 = $G110_old
$ws.Range("H111").Value2 = $H110_old
$ws.Range("I111").Value2 = $I110_old
$ws.Range("J111").Value2 = $J110_old
$ws.Range("K111").Value2 = $K110_old
$ws.Range("L111").Value2 = $L110_old
$ws.Range("M111").Value2 = $M110_old
$ws.Range("N111").Value2 = $N110_old
$ws.Range("O111").Value2 = $O110_old
$ws.Range("P111").Value2 = $P110_old
$ws.Range("Q111").Value2 = $Q110_old
$ws.Range("R111").Value2 = $R110_old
$ws.Range("S111").Value2 = $S110_old
$ws.Range("T111").Value2 = $T110_old
$ws.Range("U111").Value2 = $U110_old
$ws.Range("V111").Value2 = $V110_old
$ws.Range("W111").Value2 = $W110_old
$ws.Range("X111").Value2 = $X110_old
$ws.Range("Y111").Value2 = $Y110_old
$ws.Range("Z111").Value2 = $Z110_old
$ws.Range("AA111").Value2 = $AA110_old
$ws.Range("AB111").Value2 = $AB110_old

# --- swap row 231 <-> row 232 ---
$B231_old = $ws.Range("B231").Value2
$C231_old = $ws.Range("C231").Value2
$D231_old = $ws.Range("D231").Value2
$E231_old = $ws.Range("E231").Value2
$F231_old = $ws.Range("F231").Value2
$G231_old = $ws.Range("G231").Value2
$H231_old = $ws.Range("H231").Value2
$I231_old = $ws.Range("I231").Value2
$J231_old = $ws.Range("J231").Value2
$K231_old = $ws.Range("K231").Value2
$L231_old = $ws.Range("L231").Value2
$M231_old = $ws.Range("M231").Value2
$N231_old = $ws.Range("N231").Value2
$O231_old = $ws.Range("O231").Value2
$P231_old = $ws.Range("P231").Value2
$Q231_old = $ws.Range("Q231").Value2
$R231_old = $ws.Range("R231").Value2
$S231_old = $ws.Range("S231").Value2
$T231_old = $ws.Range("T231").Value2
$U231_old = $ws.Range("U231").Value2
$V231_old = $ws.Range("V231").Value2
$W231_old = $ws.Range("W231").Value2
$X231_old = $ws.Range("X231").Value2
$Y231_old = $ws.Range("Y231").Value2
$Z231_old = $ws.Range("Z231").Value2
$AA231_old = $ws.Range("AA231").Value2
$AB231_old = $ws.Range("AB231").Value2
$B232_old = $ws.Range("B232").Value2
$C232_old = $ws.Range("C232").Value2
$D232_old = $ws.Range("D232").Value2
$E232_old = $ws.Range("E232").Value2
$F232_old = $ws.Range("F232").Value2
$G232_old = $ws.Range("G232").Value2
$H232_old = $ws.Range("H232").Value2
$I232_old = $ws.Range("I232").Value2
$J232_old = $ws.Range("J232").Value2
$K232_old = $ws.Range("K232").Value2
$L232_old = $ws.Range("L232").Value2
$M232_old = $ws.Range("M232").Value2
$N232_old = $ws.Range("N232").Value2
$O232_old = $ws.Range("O232").Value2
$P232_old = $ws.Range("P232").Value2
$Q232_old = $ws.Range("Q232").Value2
$R232_old = $ws.Range("R232").Value2
$S232_old = $ws.Range("S232").Value2
$T232_old = $ws.Range("T232").Value2
$U232_old = $ws.Range("U232").Value2
$V232_old = $ws.Range("V232").Value2
$W232_old = $ws.Range("W232").Value2
$X232_old = $ws.Range("X232").Value2
$Y232_old = $ws.Range("Y232").Value2
$Z232_old = $ws.Range("Z232").Value2
$AA232_old = $ws.Range("AA232").Value2
$AB232_old = $ws.Range("AB232").Value2
$ws.Range("B231").Value2 = $B232_old
$ws.Range("C231").Value2 = $C232_old
$ws.Range("D231").Value2 = $D232_old
$ws.Range("E231").Value2 = $E232_old
$ws.Range("F231").Value2 = $F232_old
$ws.Range("G231").Value2 = $G232_old
$ws.Range("H231").Value2 = $H232_old
$ws.Range("I231").Value2 = $I232_old
$ws.Range("J231").Value2 = $J232_old
$ws.Range("K231").Value2 = $K232_old
$ws.Range("L231").Value2 = $L232_old
$ws.Range("M231").Value2 = $M232_old
$ws.Range("N231").Value2 = $N232_old
$ws.Range("O231").Value2 = $O232_old
$ws.Range("P231").Value2 = $P232_old
$ws.Range("Q231").Value2 = $Q232_old
$ws.Range("R231").Value2 = $R232_old
$ws.Range("S231").Value2 = $S232_old
$ws.Range("T231").Value2 = $T232_old
$ws.Range("U231").Value2 = $U232_old
$ws.Range("V231").Value2 = $V232_old
$ws.Range("W231").Value2 = $W232_old
$ws.Range("X231").Value2 = $X232_old
$ws.Range("Y231").Value2 = $Y232_old
$ws.Range("Z231").Value2 = $Z232_old
$ws.Range("AA231").Value2 = $AA232_old
$ws.Range("AB231").Value2 = $AB232_old
$ws.Range("B232").Value2 = $B231_old
$ws.Range("C232").Value2 = $C231_old
$ws.Range("D232").Value2 = $D231_old
$ws.Range("E232").Value2 = $E231_old
$ws.Range("F232").Value2 = $F231_old
$ws.Range("G232").Value2 = $G231_old
$ws.Range("H232").Value2 = $H231_old
$ws.Range("I232").Value2 = $I231_old
$ws.Range("J232").Value2 = $J231_old
$ws.Range("K232").Value2 = $K231_old
$ws.Range("L232").Value2 = $L231_old
$ws.Range("M232").Value2 = $M231_old
$ws.Range("N232").Value2 = $N231_old
$ws.Range("O232").Value2 = $O231_old
$ws.Range("P232").Value2 = $P231_old
$ws.Range("Q232").Value2 = $Q231_old
$ws.Range("R232").Value2 = $R231_old
$ws.Range("S232").Value2 = $S231_old
$ws.Range("T232").Value2 = $T231_old
$ws.Range("U232").Value2 = $U231_old
$ws.Range("V232").Value2 = $V231_old
$ws.Range("W232").Value2 = $W231_old
$ws.Range("X232").Value2 = $X231_old
$ws.Range("Y232").Value2 = $Y231_old
$ws.Range("Z232").Value2 = $Z231_old
$ws.Range("AA232").Value2 = $AA231_old
$ws.Range("AB232").Value2 = $AB231_old

# --- swap row 269 <-> row 271 ---
$B269_old = $ws.Range("B269").Value2
$C269_old = $ws.Range("C269").Value2
$D269_old = $ws.Range("D269").Value2
$E269_old = $ws.Range("E269").Value2
$F269_old = $ws.Range("F269").Value2
$G269_old = $ws.Range("G269").Value2
$H269_old = $ws.Range("H269").Value2
$I269_old = $ws.Range("I269").Value2
$J269_old = $ws.Range("J269").Value2
$K269_old = $ws.Range("K269").Value2
$L269_old = $ws.Range("L269").Value2
$M269_old = $ws.Range("M269").Value2
$N269_old = $ws.Range("N269").Value2
$O269_old = $ws.Range("O269").Value2
$P269_old = $ws.Range("P269").Value2
$Q269_old = $ws.Range("Q269").Value2
$R269_old = $ws.Range("R269").Value2
$S269_old = $ws.Range("S269").Value2
$T269_old = $ws.Range("T269").Value2
$U269_old = $ws.Range("U269").Value2
$V269_old = $ws.Range("V269").Value2
$W269_old = $ws.Range("W269").Value2
$X269_old = $ws.Range("X269").Value2
$Y269_old = $ws.Range("Y269").Value2
$Z269_old = $ws.Range("Z269").Value2
$AA269_old = $ws.Range("AA269").Value2
$AB269_old = $ws.Range("AB269").Value2
$B271_old = $ws.Range("B271").Value2
$C271_old = $ws.Range("C271").Value2
$D271_old = $ws.Range("D271").Value2
$E271_old = $ws.Range("E271").Value2
$F271_old = $ws.Range("F271").Value2
$G271_old = $ws.Range("G271").Value2
$H271_old = $ws.Range("H271").Value2
$I271_old = $ws.Range("I271").Value2
$J271_old = $ws.Range("J271").Value2
$K271_old = $ws.Range("K271").Value2
$L271_old = $ws.Range("L271").Value2
$M271_old = $ws.Range("M271").Value2
$N271_old = $ws.Range("N271").Value2
$O271_old = $ws.Range("O271").Value2
$P271_old = $ws.Range("P271").Value2
$Q271_old = $ws.Range("Q271").Value2
$R271_old = $ws.Range("R271").Value2
$S271_old = $ws.Range("S271").Value2
$T271_old = $ws.Range("T271").Value2
$U271_old = $ws.Range("U271").Value2
$V271_old = $ws.Range("V271").Value2
$W271_old = $ws.Range("W271").Value2
$X271_old = $ws.Range("X271").Value2
$Y271_old = $ws.Range("Y271").Value2
$Z271_old = $ws.Range("Z271").Value2
$AA271_old = $ws.Range("AA271").Value2
$AB271_old = $ws.Range("AB271").Value2
$ws.Range("B269").Value2 = $B271_old
$ws.Range("C269").Value2 = $C271_old
$ws.Range("D269").Value2 = $D271_old
$ws.Range("E269").Value2 = $E271_old
$ws.Range("F269").Value2 = $F271_old
$ws.Range("G269").Value2 = $G271_old
$ws.Range("H269").Value2 = $H271_old
$ws.Range("I269").Value2 = $I271_old
$ws.Range("J269").Value2 = $J271_old
$ws.Range("K269").Value2 = $K271_old
$ws.Range("L269").Value2 = $L271_old
$ws.Range("M269").Value2 = $M271_old
$ws.Range("N269").Value2 = $N271_old
$ws.Range("O269").Value2 = $O271_old
$ws.Range("P269").Value2 = $P271_old
$ws.Range("Q269").Value2 = $Q271_old
$ws.Range("R269").Value2 = $R271_old
$ws.Range("S269").Value2 = $S271_old
$ws.Range("T269").Value2 = $T271_old
$ws.Range("U269").Value2 = $U271_old
$ws.Range("V269").Value2 = $V271_old
$ws.Range("W269").Value2 = $W271_old
$ws.Range("X269").Value2 = $X271_old
$ws.Range("Y269").Value2 = $Y271_old
$ws.Range("Z269").Value2 = $Z271_old
$ws.Range("AA269").Value2 = $AA271_old
$ws.Range("AB269").Value2 = $AB271_old
$ws.Range("B271").Value2 = $B269_old
$ws.Range("C271").Value2 = $C269_old
$ws.Range("D271").Value2 = $D269_old
$ws.Range("E271").Value2 = $E269_old
$ws.Range("F271").Value2 = $F269_old
$ws.Range("G271").Value2 = $G269_old
$ws.Range("H271").Value2 = $H269_old
$ws.Range("I271").Value2 = $I269_old
$ws.Range("J271").Value2 = $J269_old
$ws.Range("K271").Value2 = $K269_old
$ws.Range("L271").Value2 = $L269_old
$ws.Range("M271").Value2 = $M269_old
$ws.Range("N271").Value2 = $N269_old
$ws.Range("O271").Value2 = $O269_old
$ws.Range("P271").Value2 = $P269_old
$ws.Range("Q271").Value2 = $Q269_old
$ws.Range("R271").Value2 = $R269_old
$ws.Range("S271").Value2 = $S269_old
$ws.Range("T271").Value2 = $T269_old
$ws.Range("U271").Value2 = $U269_old
$ws.Range("V271").Value2 = $V269_old
$ws.Range("W271").Value2 = $W269_old
$ws.Range("X271").Value2 = $X269_old
$ws.Range("Y271").Value2 = $Y269_old
$ws.Range("Z271").Value2 = $Z269_old
$ws.Range("AA271").Value2 = $AA269_old
$ws.Range("AB271").Value2 = $AB269_old
